# Add the new "2020" column (O) to the participation-rate table, mirroring
# the formatting already used by the adjacent "2019" column (N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: year 2020
$ws.Range("O4").Value = 2020
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)  # xlPasteFormats - copy N4's formatting onto O4

# Data cell: 2020 value
$ws.Range("O5").Value = 83.3
$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)  # xlPasteFormats - copy N5's formatting onto O5

# Clear the marching-ants clipboard marquee left over from the copies above
$excel.CutCopyMode = $false

# Match the saved selection state of the edited workbook
$ws.Range("O12").Select() | Out-Null
